$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.729797666666667
$ws.Range("H2").Value = 5.189393000000001
$ws.Range("I2").Value = 0.06436583050179444
$ws.Range("J2").Value = 0.06436583050179444
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 213.2505674533679
$ws.Range("R2").Value = 1919.255107080311
$ws.Range("S2").Value = 0.04017440981627448
$ws.Range("T2").Value = 0.04017440981627447
$ws.Range("G3").Value = 1.729797666666667
$ws.Range("H3").Value = 5.189393000000001
$ws.Range("I3").Value = 0.06436583050179444
$ws.Range("J3").Value = 0.06436583050179444
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 81.75033863153055
$ws.Range("R3").Value = 735.7530476837751
$ws.Range("S3").Value = 0.015400998206115
$ws.Range("T3").Value = 0.015400998206115
$ws.Range("G4").Value = 1.729797666666667
$ws.Range("H4").Value = 5.189393000000001
$ws.Range("I4").Value = 0.06436583050179444
$ws.Range("J4").Value = 0.06436583050179444
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 46.66061282444957
$ws.Range("R4").Value = 419.9455154200461
$ws.Range("S4").Value = 0.008790422479404963
$ws.Range("T4").Value = 0.008790422479404961
$ws.Range("I5").Value = 0.2200595722726403
$ws.Range("J5").Value = 0.2200595722726403
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 729.0798284561509
$ws.Range("R5").Value = 6561.718456105358
$ws.Range("S5").Value = 0.1373518118472604
$ws.Range("T5").Value = 0.1373518118472604
$ws.Range("I6").Value = 0.2200595722726403
$ws.Range("J6").Value = 0.2200595722726403
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("S6").Value = 0.05265428957239792
$ws.Range("T6").Value = 0.05265428957239791
$ws.Range("I7").Value = 0.2200595722726403
$ws.Range("J7").Value = 0.2200595722726403
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("S7").Value = 0.03005347085298199
$ws.Range("T7").Value = 0.03005347085298199
$ws.Range("H8").Value = 57.69206699999999
$ws.Range("I8").Value = 0.7155745972255653
$ws.Range("J8").Value = 0.7155745972255653
$ws.Range("M8").Value = 123.2806423333333
$ws.Range("N8").Value = 369.841927
$ws.Range("O8").Value = 0.6241574062367528
$ws.Range("P8").Value = 0.6241574062367526
$ws.Range("Q8").Value = 2370.771692432567
$ws.Range("R8").Value = 21336.94523189311
$ws.Range("S8").Value = 0.4466311845732179
$ws.Range("T8").Value = 0.4466311845732178
$ws.Range("H9").Value = 57.69206699999999
$ws.Range("I9").Value = 0.7155745972255653
$ws.Range("J9").Value = 0.7155745972255653
$ws.Range("O9").Value = 0.2392728888301323
$ws.Range("P9").Value = 0.2392728888301322
$ws.Range("Q9").Value = 908.8434839301915
$ws.Range("R9").Value = 8179.591355371724
$ws.Range("S9").Value = 0.1712176010516193
$ws.Range("T9").Value = 0.1712176010516193
$ws.Range("H10").Value = 57.69206699999999
$ws.Range("I10").Value = 0.7155745972255653
$ws.Range("J10").Value = 0.7155745972255653
$ws.Range("O10").Value = 0.136569704933115
$ws.Range("P10").Value = 0.136569704933115
$ws.Range("Q10").Value = 518.7402845244526
$ws.Range("R10").Value = 4668.662560720074
$ws.Range("S10").Value = 0.09772581160072809
$ws.Range("T10").Value = 0.09772581160072807
